$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values that changed ---
$ws.Range("A5").Value = "branch_of_study"
$ws.Range("A7").Value = "job/intern/learning/corporate"

# --- Add new column B attribute values for rows 2-5 ---
$ws.Range("B2").Value = "times_executed"
$ws.Range("B3").Value = "error_count"
$ws.Range("B4").Value = "time_taken_to_solve"
$ws.Range("B5").Value = "keystrokes (backspace)"

# --- Append new attribute rows (10-16) ---
$ws.Range("A10").Value = "branch_of_study"
$ws.Range("A11").Value = "math/non_math"
$ws.Range("A12").Value = "working/non-working"
$ws.Range("A13").Value = "total_login_time"
$ws.Range("A14").Value = "total_classes_attended"
$ws.Range("A15").Value = "total_assignments_completed"
$ws.Range("A16").Value = "free_time_after_college"

# --- Resize columns A:B to fit the new (longer) content ---
$ws.Columns("A:B").AutoFit()

# --- Move the selection cursor like in the saved file ---
$ws.Range("E8").Select()
